$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Edit the "giangVien" row description (B7): insert "batDauTu, ketThucNgay, "
# right after "ngaySinh, " and before "soDienThoai" in the rich-text cell.
$cell = $ws.Range("B7")

$run1 = "giangVienId,"
$run2 = "  firstName, lastName, gioiTinh, ngaySinh, batDauTu, ketThucNgay, soDienThoai, gmail, diaChi, "
$run3 = "idCoSo , taiKhoanId"
$newText = $run1 + $run2 + $run3

$cell.Value = $newText

$run1Len = $run1.Length
$run2Len = $run2.Length
$run3Len = $run3.Length
$totalLen = $newText.Length

# Re-apply the original per-run colors (they get reset to default when the
# value is overwritten): run1 = red, run2 = black, run3 = blue (accent1).
$c1 = $cell.Characters(1, $run1Len)
$c1.Font.Color = 255

$c2 = $cell.Characters($run1Len + 1, $run2Len)
$c2.Font.Color = 0

$c3 = $cell.Characters($totalLen - $run3Len + 1, $run3Len)
$c3.Font.Color = 0x5B + (0x9B * 256) + (0xD5 * 65536)

# --- Move the active selection to B3 (cosmetic, matches the saved view state)
$ws.Range("B3").Select() | Out-Null
